# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Each entry below maps a row number to its new F-column value for a given sheet.

$wb = $excel.ActiveWorkbook

$sheet1Updates = @{
    2  = 21
    4  = 85
    5  = 20
    6  = 550
    7  = 1702
    10 = 32
    11 = 1722
    13 = 95
    14 = 415
    16 = 196
    17 = 14
    21 = 495
    24 = 234
    25 = 253
}

$sheet4Updates = @{
    2  = 21
    4  = 85
    5  = 20
    6  = 550
    7  = 1702
    11 = 32
    12 = 1722
    14 = 95
    15 = 415
    16 = 268
    17 = 196
    18 = 14
    22 = 495
    25 = 234
    26 = 253
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
